# Update simulated-game transition-probability matrix on Sheet1 with
# refreshed values following more simulated games / reworked sim logic.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.208984375
$ws.Range("C2").Value = 0.5078125
$ws.Range("J2").Value = 0.0234375
$ws.Range("P2").Value = 0.150390625
$ws.Range("S2").Value = 0.109375
$ws.Range("B3").Value = 0.01486988847583643
$ws.Range("C3").Value = 0.02973977695167286
$ws.Range("J3").Value = 0.04089219330855019
$ws.Range("P3").Value = 0.7323420074349443
$ws.Range("S3").Value = 0.1821561338289963
$ws.Range("J4").Value = 0.07954545454545454
$ws.Range("O4").Value = 0.01136363636363636
$ws.Range("P4").Value = 0.6363636363636364
$ws.Range("S4").Value = 0.2727272727272727
$ws.Range("P5").Value = 0.6
$ws.Range("S5").Value = 0.4
$ws.Range("B6").Value = 0.0457516339869281
$ws.Range("D6").Value = 0.01525054466230937
$ws.Range("F6").Value = 0.07407407407407407
$ws.Range("J6").Value = 0.2549019607843137
$ws.Range("O6").Value = 0.03267973856209151
$ws.Range("Q6").Value = 0.1546840958605664
$ws.Range("R6").Value = 0.07625272331154684
$ws.Range("S6").Value = 0.3464052287581699
$ws.Range("B7").Value = 0.09315068493150686
$ws.Range("D7").Value = 0.02465753424657534
$ws.Range("E7").Value = 0.005479452054794521
$ws.Range("F7").Value = 0.06027397260273973
$ws.Range("J7").Value = 0.1342465753424658
$ws.Range("O7").Value = 0.0273972602739726
$ws.Range("Q7").Value = 0.1753424657534247
$ws.Range("R7").Value = 0.0821917808219178
$ws.Range("S7").Value = 0.3972602739726027
$ws.Range("B8").Value = 0.08668341708542714
$ws.Range("D8").Value = 0.01633165829145729
$ws.Range("E8").Value = 0.001256281407035176
$ws.Range("F8").Value = 0.06407035175879397
$ws.Range("J8").Value = 0.1218592964824121
$ws.Range("O8").Value = 0.02010050251256281
$ws.Range("Q8").Value = 0.1947236180904523
$ws.Range("R8").Value = 0.10678391959799
$ws.Range("S8").Value = 0.3881909547738693
$ws.Range("B9").Value = 0.0966183574879227
$ws.Range("D9").Value = 0.00966183574879227
$ws.Range("E9").Value = 0.002415458937198068
$ws.Range("F9").Value = 0.08695652173913043
$ws.Range("J9").Value = 0.1280193236714976
$ws.Range("O9").Value = 0.04830917874396135
$ws.Range("Q9").Value = 0.1714975845410628
$ws.Range("R9").Value = 0.1328502415458937
$ws.Range("S9").Value = 0.3236714975845411
$ws.Range("B10").Value = 0.08738980452280568
$ws.Range("D10").Value = 0.02069758528171713
$ws.Range("E10").Value = 0.0007665772326561902
$ws.Range("F10").Value = 0.0697585281717133
$ws.Range("J10").Value = 0.1215024913760061
$ws.Range("O10").Value = 0.02453047144499809
$ws.Range("Q10").Value = 0.2108087389804523
$ws.Range("R10").Value = 0.1061709467228823
$ws.Range("S10").Value = 0.3583748562667689
$ws.Range("G11").Value = 0.1393728222996516
$ws.Range("J11").Value = 0.1045296167247387
$ws.Range("K11").Value = 0.1846689895470383
$ws.Range("L11").Value = 0.5592334494773519
$ws.Range("S11").Value = 0.01219512195121951
$ws.Range("G12").Value = 0.7236024844720497
$ws.Range("J12").Value = 0.2546583850931677
$ws.Range("K12").Value = 0.006211180124223602
$ws.Range("L12").Value = 0.003105590062111801
$ws.Range("S12").Value = 0.0124223602484472
$ws.Range("G13").Value = 0.6956521739130435
$ws.Range("J13").Value = 0.2608695652173913
$ws.Range("S13").Value = 0.04347826086956522
$ws.Range("F15").Value = 0.02188183807439825
$ws.Range("H15").Value = 0.1444201312910285
$ws.Range("I15").Value = 0.06345733041575492
$ws.Range("J15").Value = 0.3282275711159737
$ws.Range("K15").Value = 0.06345733041575492
$ws.Range("M15").Value = 0.01969365426695843
$ws.Range("O15").Value = 0.0700218818380744
$ws.Range("S15").Value = 0.2888402625820569
$ws.Range("F16").Value = 0.02760736196319018
$ws.Range("H16").Value = 0.1564417177914111
$ws.Range("I16").Value = 0.0705521472392638
$ws.Range("J16").Value = 0.4325153374233129
$ws.Range("K16").Value = 0.1257668711656442
$ws.Range("M16").Value = 0.02760736196319018
$ws.Range("O16").Value = 0.03680981595092025
$ws.Range("S16").Value = 0.1226993865030675
$ws.Range("F17").Value = 0.01779755283648498
$ws.Range("H17").Value = 0.1635150166852058
$ws.Range("I17").Value = 0.1078976640711902
$ws.Range("J17").Value = 0.4271412680756396
$ws.Range("K17").Value = 0.09232480533926585
$ws.Range("M17").Value = 0.02669632925472748
$ws.Range("N17").Value = 0.001112347052280311
$ws.Range("O17").Value = 0.06451612903225806
$ws.Range("S17").Value = 0.09899888765294772
$ws.Range("F18").Value = 0.02702702702702703
$ws.Range("H18").Value = 0.1995841995841996
$ws.Range("I18").Value = 0.103950103950104
$ws.Range("J18").Value = 0.4074844074844075
$ws.Range("K18").Value = 0.07276507276507277
$ws.Range("M18").Value = 0.01871101871101871
$ws.Range("O18").Value = 0.05405405405405406
$ws.Range("S18").Value = 0.1164241164241164
$ws.Range("F19").Value = 0.01455301455301455
$ws.Range("H19").Value = 0.1808731808731809
$ws.Range("I19").Value = 0.09022869022869023
$ws.Range("J19").Value = 0.3995841995841996
$ws.Range("K19").Value = 0.1081081081081081
$ws.Range("M19").Value = 0.01954261954261954
$ws.Range("N19").Value = 0.0004158004158004158
$ws.Range("O19").Value = 0.0700218818380744
$ws.Range("S19").Value = 0.1239085239085239
